$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 810, shifting existing rows 810..874 down to 811..875
$ws.Rows.Item(810).Insert()

# Populate the newly inserted row 810 with the new record
$ws.Cells.Item(810,1).Value = 3
$ws.Cells.Item(810,2).Value = "Femacal de La Calera"
$ws.Cells.Item(810,3).Value = "Coquimbo"
$ws.Cells.Item(810,4).Value = 45013
$ws.Cells.Item(810,5).Value = 5
$ws.Cells.Item(810,6).Value = 100112045
$ws.Cells.Item(810,7).Value = "Zapallo"
$ws.Cells.Item(810,8).Value = "Paine"
$ws.Cells.Item(810,9).Value = "1a (cosecha)"
$ws.Cells.Item(810,10).Value = 190
$ws.Cells.Item(810,11).Value = 550
$ws.Cells.Item(810,12).Value = 550
$ws.Cells.Item(810,13).Value = 550
$ws.Cells.Item(810,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(810,15).Value = "Provincia de Talca"
$ws.Cells.Item(810,16).Value = 550
$ws.Cells.Item(810,17).Value = 1
$ws.Cells.Item(810,18).Value = "Hortaliza"
